# Add files via upload
# The underlying data edit in this revision is the addition of "CU"
# markers in the CRUD matrix for the "Pull Inventory Report" (row 12)
# and "Pull-Inventory Report Definition" (row 13) rows, in the
# PULL-INVENTORY-LINE (column I) and PULL-INVENTORY (column J) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I12").Value = "CU"
$ws.Range("J12").Value = "CU"
$ws.Range("I13").Value = "CU"
$ws.Range("J13").Value = "CU"

# Selection moved to M1 in the saved file.
[void]$ws.Range("M1").Select()
